{"js": "// Locate the paragraph that ends the previous results block:\n// \"1000 val\u00f3s k\u00e9p eset\u00e9n: 41,25%-os pontoss\u00e1g\" (the last occurrence of the\n// \"%-os pontoss\u00e1g\" results list), then append the new narrative + results\n// paragraphs directly after it, in document order.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"1000 val\u00f3s k\u00e9p eset\u00e9n\") !== -1) {\n    anchor = p;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Anchor paragraph ('1000 val\u00f3s k\u00e9p eset\u00e9n...') not found\");\n}\n\n// --- new paragraph 1/15 (empty) ---\nanchor = anchor.insertParagraph(\"\", \"After\");\n\n// --- new paragraph 2/15  ---\nanchor = anchor.insertParagraph(\"Ezut\u00e1n tov\u00e1bbi val\u00f3s k\u00e9pek hozz\u00e1ad\u00e1s\u00e1val pr\u00f3b\u00e1ltam n\u00f6velni a val\u00f3s k\u00e9pet felismer\u00e9s\u00e9nek \", \"After\");\nanchor.getRange(\"End\").insertText(\"hat\u00e9konys\u00e1g\u00e1t. Ez azonban nem jav\u00edtott a teljes\u00edtm\u00e9nyen.\", \"End\");\n\n// --- new paragraph 3/15 (empty) ---\nanchor = anchor.insertParagraph(\"\", \"After\");\n\n// --- new paragraph 4/15  ---\nanchor = anchor.insertParagraph(\"Majd egy Color Jitter f\u00fcggv\u00e9nyt pr\u00f3b\u00e1ltam a k\u00e9pekre rakni. Ez v\u00e9letlenszer\u0171en v\u00e1ltoztatja a \", \"After\");\nanchor.getRange(\"End\").insertText(\"k\u00e9pek sz\u00edn\u00e9t egy tartom\u00e1nyban\", \"End\");\nanchor.getRange(\"End\").insertText(\", ami seg\u00edt \u00e1ltal\u00e1nosabb\u00e1 tenni a neur\u00e1lis h\u00e1l\u00f3t, mivel nem tanul meg annyira specifikus mint\u00e1kat (t\u00faltanul\u00e1s ellen hat\u00e1sos). Sajnos ez sem hozott v\u00e1ltoz\u00e1st.\", \"End\");\n\n// --- new paragraph 5/15 (empty) ---\nanchor = anchor.insertParagraph(\"\", \"After\");\n\n// --- new paragraph 6/15  ---\nanchor = anchor.insertParagraph(\"A modell kimenete k\u00e9t sz\u00e1m, amelyb\u0151l j\u00f3l lehet sz\u00e1zal\u00e9kosan sz\u00e1molni, hogy mekkora es\u00e9llyel AI gener\u00e1lt vagy val\u00f3s a k\u00e9p. Mivel tov\u00e1bbra is az AI gener\u00e1lt k\u00e9pek eset\u00e9n magas a pontoss\u00e1g, a val\u00f3s k\u00e9pekn\u00e9l viszont nem annyira, \", \"After\");\nanchor.getRange(\"End\").insertText(\"\u00e9s l\u00e1tszik, hogy az AI gener\u00e1lt eset\u00e9n nagy magabiztoss\u00e1ggal \u00edt\u00e9li AI gener\u00e1ltnak, m\u00edg a val\u00f3sn\u00e1l a k\u00e9t sz\u00e1m k\u00f6zelebb van egym\u00e1shoz (m\u00e9g ha az AI gener\u00e1lt tov\u00e1bbra is magasabb), \", \"End\");\nanchor.getRange(\"End\").insertText(\"\u00edgy \", \"End\");\nanchor.getRange(\"End\").insertText(\"a sz\u00e1zal\u00e9kos \u00e9rt\u00e9kek eltol\u00e1s\u00e1val tov\u00e1bb jav\u00edthat\u00f3 a pontoss\u00e1g. \", \"End\");\nanchor.getRange(\"End\").insertText(\"Az eredm\u00e9nyek:\", \"End\");\n\n// --- new paragraph 7/15 (empty) ---\nanchor = anchor.insertParagraph(\"\", \"After\");\n\n// --- new paragraph 8/15  ---\nanchor = anchor.insertParagraph(\"50% eset\u00e9n (alapeset): 65,83%\", \"After\");\n\n// --- new paragraph 9/15  ---\nanchor = anchor.insertParagraph(\"60% eset\u00e9n: \", \"After\");\nanchor.getRange(\"End\").insertText(\"69,28%\", \"End\");\n\n// --- new paragraph 10/15  ---\nanchor = anchor.insertParagraph(\"70% eset\u00e9n:\", \"After\");\nanchor.getRange(\"End\").insertText(\" 71,16%\", \"End\");\n\n// --- new paragraph 11/15  ---\nanchor = anchor.insertParagraph(\"80% eset\u00e9n: \", \"After\");\nanchor.getRange(\"End\").insertText(\"76,49%\", \"End\");\n\n// --- new paragraph 12/15  ---\nanchor = anchor.insertParagraph(\"90% eset\u00e9n: \", \"After\");\nanchor.getRange(\"End\").insertText(\"82,13%\", \"End\");\n\n// --- new paragraph 13/15 (empty) ---\nanchor = anchor.insertParagraph(\"\", \"After\");\n\n// --- new paragraph 14/15  ---\nanchor = anchor.insertParagraph(\"Teh\u00e1t azzal, hogy csak akkor \u00edt\u00e9l AI gener\u00e1ltnak egy k\u00e9pet, hogyha 0,9-n\u00e9l nagyobb az els\u0151 \u00e9rt\u00e9k \u00e9s 0,1-n\u00e9l kisebb a m\u00e1sodik, a pontoss\u00e1g 16%-ot javult.\", \"After\");\n\n// --- new paragraph 15/15  ---\nanchor = anchor.insertParagraph(\"\u00cdgy a val\u00f3s k\u00e9pekn\u00e9l a pontoss\u00e1g 66,25%, m\u00edg az AI gener\u00e1ltn\u00e1l 98,14%.\", \"After\");\n\nawait context.sync();", "ps1": "# Locate the paragraph that ends the previous results block:\n# \"1000 val\u00f3s k\u00e9p eset\u00e9n: 41,25%-os pontoss\u00e1g\" (the last occurrence of the\n# \"%-os pontoss\u00e1g\" results list), then append the new narrative + results\n# paragraphs directly after it, in document order.\n\nfunction New-ParaAfter($rng) {\n    # Word COM: InsertParagraphAfter() splices in a paragraph mark right\n    # after $rng but does NOT move $rng itself, so step past the inserted\n    # mark (+1 char) to land inside the freshly created paragraph.\n    $rng.InsertParagraphAfter()\n    $rng.Start = $rng.End + 1\n    $rng.End = $rng.Start\n    return $rng\n}\n\nfunction Append-Run($rng, [string]$text) {\n    $rng.InsertAfter($text)\n    $rng.Start = $rng.End\n    return $rng\n}\n\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"1000 val\u00f3s k\u00e9p eset\u00e9n: 41,25%-os pontoss\u00e1g\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Anchor paragraph ('1000 val\u00f3s k\u00e9p eset\u00e9n...') not found\"\n}\n$rng = $find.Parent\n$rng.Collapse(0)  # wdCollapseEnd\n\n# --- new paragraph 1/15 (empty) ---\n$rng = New-ParaAfter $rng\n\n# --- new paragraph 2/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"Ezut\u00e1n tov\u00e1bbi val\u00f3s k\u00e9pek hozz\u00e1ad\u00e1s\u00e1val pr\u00f3b\u00e1ltam n\u00f6velni a val\u00f3s k\u00e9pet felismer\u00e9s\u00e9nek \"\n$rng = Append-Run $rng \"hat\u00e9konys\u00e1g\u00e1t. Ez azonban nem jav\u00edtott a teljes\u00edtm\u00e9nyen.\"\n\n# --- new paragraph 3/15 (empty) ---\n$rng = New-ParaAfter $rng\n\n# --- new paragraph 4/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"Majd egy Color Jitter f\u00fcggv\u00e9nyt pr\u00f3b\u00e1ltam a k\u00e9pekre rakni. Ez v\u00e9letlenszer\u0171en v\u00e1ltoztatja a \"\n$rng = Append-Run $rng \"k\u00e9pek sz\u00edn\u00e9t egy tartom\u00e1nyban\"\n$rng = Append-Run $rng \", ami seg\u00edt \u00e1ltal\u00e1nosabb\u00e1 tenni a neur\u00e1lis h\u00e1l\u00f3t, mivel nem tanul meg annyira specifikus mint\u00e1kat (t\u00faltanul\u00e1s ellen hat\u00e1sos). Sajnos ez sem hozott v\u00e1ltoz\u00e1st.\"\n\n# --- new paragraph 5/15 (empty) ---\n$rng = New-ParaAfter $rng\n\n# --- new paragraph 6/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"A modell kimenete k\u00e9t sz\u00e1m, amelyb\u0151l j\u00f3l lehet sz\u00e1zal\u00e9kosan sz\u00e1molni, hogy mekkora es\u00e9llyel AI gener\u00e1lt vagy val\u00f3s a k\u00e9p. Mivel tov\u00e1bbra is az AI gener\u00e1lt k\u00e9pek eset\u00e9n magas a pontoss\u00e1g, a val\u00f3s k\u00e9pekn\u00e9l viszont nem annyira, \"\n$rng = Append-Run $rng \"\u00e9s l\u00e1tszik, hogy az AI gener\u00e1lt eset\u00e9n nagy magabiztoss\u00e1ggal \u00edt\u00e9li AI gener\u00e1ltnak, m\u00edg a val\u00f3sn\u00e1l a k\u00e9t sz\u00e1m k\u00f6zelebb van egym\u00e1shoz (m\u00e9g ha az AI gener\u00e1lt tov\u00e1bbra is magasabb), \"\n$rng = Append-Run $rng \"\u00edgy \"\n$rng = Append-Run $rng \"a sz\u00e1zal\u00e9kos \u00e9rt\u00e9kek eltol\u00e1s\u00e1val tov\u00e1bb jav\u00edthat\u00f3 a pontoss\u00e1g. \"\n$rng = Append-Run $rng \"Az eredm\u00e9nyek:\"\n\n# --- new paragraph 7/15 (empty) ---\n$rng = New-ParaAfter $rng\n\n# --- new paragraph 8/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"50% eset\u00e9n (alapeset): 65,83%\"\n\n# --- new paragraph 9/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"60% eset\u00e9n: \"\n$rng = Append-Run $rng \"69,28%\"\n\n# --- new paragraph 10/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"70% eset\u00e9n:\"\n$rng = Append-Run $rng \" 71,16%\"\n\n# --- new paragraph 11/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"80% eset\u00e9n: \"\n$rng = Append-Run $rng \"76,49%\"\n\n# --- new paragraph 12/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"90% eset\u00e9n: \"\n$rng = Append-Run $rng \"82,13%\"\n\n# --- new paragraph 13/15 (empty) ---\n$rng = New-ParaAfter $rng\n\n# --- new paragraph 14/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"Teh\u00e1t azzal, hogy csak akkor \u00edt\u00e9l AI gener\u00e1ltnak egy k\u00e9pet, hogyha 0,9-n\u00e9l nagyobb az els\u0151 \u00e9rt\u00e9k \u00e9s 0,1-n\u00e9l kisebb a m\u00e1sodik, a pontoss\u00e1g 16%-ot javult.\"\n\n# --- new paragraph 15/15  ---\n$rng = New-ParaAfter $rng\n$rng = Append-Run $rng \"\u00cdgy a val\u00f3s k\u00e9pekn\u00e9l a pontoss\u00e1g 66,25%, m\u00edg az AI gener\u00e1ltn\u00e1l 98,14%.\"\n"}
